$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 989.6667
$ws.Range("J53").Value = 4050
$ws.Range("L53").Value = 4050
$ws.Range("N53").Value = -5324

$ws.Range("H62").Value = 9633.333
$ws.Range("I62").Value = 9750
$ws.Range("J62").Value = 9400
$ws.Range("K62").Value = 9750
$ws.Range("L62").Value = 9400
$ws.Range("M62").Value = -9126
$ws.Range("N62").Value = -10648

$ws.Range("H65").Value = 9633.333
$ws.Range("I65").Value = 9750
$ws.Range("J65").Value = 9400
$ws.Range("K65").Value = 48750
$ws.Range("L65").Value = 47000
$ws.Range("M65").Value = -45630
$ws.Range("N65").Value = -53240

$ws.Range("H98").Value = 471.6
$ws.Range("I98").Value = 489.5
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 489.5
$ws.Range("L98").Value = 400
$ws.Range("M98").Value = 1008.5
$ws.Range("N98").Value = -3396

$ws.Range("H122").Value = 471.6
$ws.Range("I122").Value = 489.5
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 1468.5
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = 981.5
$ws.Range("N122").Value = -6100

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H141").Value = 1333.3334
$ws.Range("I141").Value = 1333.3334
$ws.Range("K141").Value = 4000.0002
$ws.Range("M141").Value = 1179.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 744.4286
$ws.Range("I2").Value = 642.2
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 642.2
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -529.2
$ws.Range("N2").Value = -1226

$ws.Range("H45").Value = 2613.4092
$ws.Range("I45").Value = 2330
$ws.Range("K45").Value = 2330
$ws.Range("M45").Value = -1953

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H116").Value = 744.4286
$ws.Range("I116").Value = 642.2
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 642.2
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 1651.8
$ws.Range("N116").Value = -5588

$ws.Range("H128").Value = 299999
$ws.Range("I128").Value = 299999
$ws.Range("K128").Value = 299999
$ws.Range("M128").Value = -295019

$ws.Range("H132").Value = 2560.3333
$ws.Range("I132").Value = 2089.6
$ws.Range("J132").Value = 4914
$ws.Range("K132").Value = 6268.799999999999
$ws.Range("L132").Value = 14742
$ws.Range("M132").Value = -3738.799999999999
$ws.Range("N132").Value = -19802

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 744.4286
$ws.Range("I3").Value = 642.2
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 642.2
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -528.2
$ws.Range("N3").Value = -1228

$ws.Range("H99").Value = 7124.75
$ws.Range("I99").Value = 7124.75
$ws.Range("K99").Value = 7124.75
$ws.Range("M99").Value = -5626.75

$ws.Range("H134").Value = 5023.25
$ws.Range("I134").Value = 2697.6667
$ws.Range("J134").Value = 12000
$ws.Range("K134").Value = 8093.000100000001
$ws.Range("L134").Value = 36000
$ws.Range("M134").Value = -5558.000100000001
$ws.Range("N134").Value = -41070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19298.4
$ws.Range("I31").Value = 19298.4
$ws.Range("K31").Value = 19298.4
$ws.Range("M31").Value = -19003.4

$ws.Range("H34").Value = 19298.4
$ws.Range("I34").Value = 19298.4
$ws.Range("K34").Value = 19298.4
$ws.Range("M34").Value = -19096.4

$ws.Range("H134").Value = 1687.5
$ws.Range("I134").Value = 1687.5
$ws.Range("K134").Value = 5062.5
$ws.Range("M134").Value = -2527.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 900
$ws.Range("J50").Value = 900
$ws.Range("L50").Value = 2700
$ws.Range("N50").Value = -3662

$ws.Range("H53").Value = 900
$ws.Range("J53").Value = 900
$ws.Range("L53").Value = 2700
$ws.Range("N53").Value = -3662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 674.44446
$ws.Range("I107").Value = 379.33334
$ws.Range("J107").Value = 1264.6666
$ws.Range("K107").Value = 379.33334
$ws.Range("L107").Value = 1264.6666
$ws.Range("M107").Value = 1540.66666
$ws.Range("N107").Value = -5104.6666

$ws.Range("H126").Value = 2080
$ws.Range("I126").Value = 2120
$ws.Range("K126").Value = 6360
$ws.Range("M126").Value = -3890

$ws.Range("H132").Value = 5408.5
$ws.Range("I132").Value = 5408.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16225.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -13695.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5999.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 5999.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 5999.5
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -6223.5

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H82").Value = 1999.4286
$ws.Range("I82").Value = 1999.4286
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1999.4286
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1638.4286
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 1999.4286
$ws.Range("I85").Value = 1999.4286
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1999.4286
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -751.4286
$ws.Range("N85").ClearContents()

$ws.Range("H122").Value = 5393.25
$ws.Range("J122").Value = 4000
$ws.Range("L122").Value = 12000
$ws.Range("N122").Value = -16900

$ws.Range("H126").Value = 5999.5
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5999.5
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 17998.5
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -22938.5

$ws.Range("H132").Value = 30000
$ws.Range("I132").Value = 30833.334
$ws.Range("K132").Value = 92500.00199999999
$ws.Range("M132").Value = -89970.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H132").Value = 5151.3335
$ws.Range("I132").Value = 5151.3335
$ws.Range("K132").Value = 15454.0005
$ws.Range("M132").Value = -12924.0005
